$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header / label columns F:I, mirroring A:D (shared strings 3,4,5,6 and 0/1/2) ---
$ws.Range("F1").Value = "5m"
$ws.Range("G1").Value = "30m"
$ws.Range("H1").Value = "1H"
$ws.Range("I1").Value = "4H"

$ws.Range("F2").Value = "Khoi luong"
$ws.Range("G2").Value = "Khoi luong"
$ws.Range("H2").Value = "Khoi luong"
$ws.Range("I2").Value = "Khoi luong"

$ws.Range("F5").Value = "RSI"
$ws.Range("G5").Value = "RSI"
$ws.Range("H5").Value = "RSI"
$ws.Range("I5").Value = "RSI"

$ws.Range("F8").Value = [char]0x2205
$ws.Range("G8").Value = [char]0x2205
$ws.Range("H8").Value = [char]0x2205
$ws.Range("I8").Value = [char]0x2205

$ws.Range("F9").Value = [char]0x2205
$ws.Range("G9").Value = [char]0x2205
$ws.Range("H9").Value = [char]0x2205
$ws.Range("I9").Value = [char]0x2205

# --- Updated values for existing A:D columns (rows 3,4,6,7) ---
$ws.Range("A3").Value = 70
$ws.Range("B3").Value = 333
$ws.Range("C3").Value = 333
$ws.Range("D3").Value = 4602

$ws.Range("A4").Value = 353.85
$ws.Range("B4").Value = 1720.2
$ws.Range("C4").Value = 2904.25
$ws.Range("D4").Value = 12999.9

$ws.Range("A6").Value = 70.03
$ws.Range("B6").Value = 53.82
$ws.Range("C6").Value = 55.22
$ws.Range("D6").Value = 62.11

$ws.Range("A7").Value = 54.78
$ws.Range("B7").Value = 47.34
$ws.Range("C7").Value = 58.49
$ws.Range("D7").Value = 66.51

# --- New values for F:I columns (rows 3,4,6,7) ---
$ws.Range("F3").Value = 366
$ws.Range("G3").Value = 821
$ws.Range("H3").Value = 821
$ws.Range("I3").Value = 821

$ws.Range("F4").Value = 412.55
$ws.Range("G4").Value = 1634.3
$ws.Range("H4").Value = 3135.15
$ws.Range("I4").Value = 12810.9

$ws.Range("F6").Value = 43.64
$ws.Range("G6").Value = 36.52
$ws.Range("H6").Value = 43.01
$ws.Range("I6").Value = 57.78

$ws.Range("F7").Value = 41.16
$ws.Range("G7").Value = 48.64
$ws.Range("H7").Value = 58.28
$ws.Range("I7").Value = 66.2

# --- Row 39-42: G column becomes formulas referencing the new F:I data ---
$ws.Range("G39").Formula = "=A4-F4"
$ws.Range("G40").Formula = "=B4-G4"
$ws.Range("G41").Formula = "=C4-H4"
$ws.Range("G42").Formula = "=D4-I4"

# --- Column widths for the new columns F:I (best effort; engine quantizes to 1/6 char steps) ---
$ws.Columns.Item(6).ColumnWidth = 9.592447916666666
$ws.Columns.Item(7).ColumnWidth = 9.736979166666666
$ws.Columns.Item(8).ColumnWidth = 9.736979166666666
$ws.Columns.Item(9).ColumnWidth = 9.877604166666666

# --- Sheet view: scroll position + selection ---
$ws.Range("I36").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
